$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2376
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11880
$ws.Range("N65").ClearContents()

$ws.Range("H70").Value = 2312.5715
$ws.Range("I70").Value = 1887.6
$ws.Range("J70").Value = 3375
$ws.Range("K70").Value = 5662.799999999999
$ws.Range("L70").Value = 10125
$ws.Range("M70").Value = -5392.799999999999
$ws.Range("N70").Value = -10665

$ws.Range("H73").Value = 2312.5715
$ws.Range("I73").Value = 1887.6
$ws.Range("J73").Value = 3375
$ws.Range("K73").Value = 5662.799999999999
$ws.Range("L73").Value = 10125
$ws.Range("M73").Value = -4726.799999999999
$ws.Range("N73").Value = -11997

$ws.Range("H80").Value = 1226
$ws.Range("I80").Value = 800
$ws.Range("J80").Value = 1439
$ws.Range("K80").Value = 2400
$ws.Range("L80").Value = 4317
$ws.Range("M80").Value = -1402
$ws.Range("N80").Value = -6313

$ws.Range("H83").Value = 1226
$ws.Range("I83").Value = 800
$ws.Range("J83").Value = 1439
$ws.Range("K83").Value = 7200
$ws.Range("L83").Value = 12951
$ws.Range("M83").Value = -2208
$ws.Range("N83").Value = -22935

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws.Range("H111").Value = 2016.25
$ws.Range("I111").Value = 2439.1
$ws.Range("J111").Value = 1311.5
$ws.Range("K111").Value = 7317.299999999999
$ws.Range("L111").Value = 3934.5
$ws.Range("M111").Value = -4250.299999999999
$ws.Range("N111").Value = -10068.5

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -8508

$ws.Range("H132").Value = 96076.14
$ws.Range("I132").Value = 111418.164
$ws.Range("J132").Value = 4024
$ws.Range("K132").Value = 334254.492
$ws.Range("L132").Value = 12072
$ws.Range("M132").Value = -331724.492
$ws.Range("N132").Value = -17132

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 1480.6666
$ws.Range("I39").Value = 1480.6666
$ws.Range("K39").Value = 1480.6666
$ws.Range("M39").Value = -960.6666

$ws.Range("H139").Value = 59998.332
$ws.Range("J139").Value = 59998.332
$ws.Range("L139").Value = 59998.332
$ws.Range("N139").Value = -70278.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 38854.75
$ws.Range("I26").Value = 30139.666
$ws.Range("K26").Value = 30139.666
$ws.Range("M26").Value = -29847.666

$ws.Range("H29").Value = 1865.5883
$ws.Range("J29").Value = 1866.6666
$ws.Range("L29").Value = 1866.6666
$ws.Range("N29").Value = -2444.6666

$ws.Range("H105").Value = 1594.5
$ws.Range("I105").Value = 1594.5
$ws.Range("K105").Value = 1594.5
$ws.Range("M105").Value = 152.5

$ws.Range("H106").Value = 24187.4
$ws.Range("J106").Value = 24187.4
$ws.Range("L106").Value = 24187.4
$ws.Range("N106").Value = -26711.4

$ws.Range("H107").Value = 41329.5
$ws.Range("I107").Value = 51276.375
$ws.Range("J107").Value = 1542
$ws.Range("K107").Value = 51276.375
$ws.Range("L107").Value = 1542
$ws.Range("M107").Value = -49356.375
$ws.Range("N107").Value = -5382

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 175.88889
$ws.Range("I7").Value = 191.95238
$ws.Range("J7").Value = 119.666664
$ws.Range("K7").Value = 191.95238
$ws.Range("L7").Value = 119.666664
$ws.Range("M7").Value = -78.95238
$ws.Range("N7").Value = -345.666664

$ws.Range("H22").Value = 615
$ws.Range("I22").Value = 646.1429
$ws.Range("K22").Value = 646.1429
$ws.Range("M22").Value = -296.1429000000001

$ws.Range("H51").Value = 33455.715
$ws.Range("I51").Value = 7500
$ws.Range("J51").Value = 37781.668
$ws.Range("K51").Value = 7500
$ws.Range("L51").Value = 37781.668
$ws.Range("M51").Value = -6764
$ws.Range("N51").Value = -39253.668

$ws.Range("H61").Value = 33455.715
$ws.Range("I61").Value = 7500
$ws.Range("J61").Value = 37781.668
$ws.Range("K61").Value = 7500
$ws.Range("L61").Value = 37781.668
$ws.Range("M61").Value = -7152
$ws.Range("N61").Value = -38477.668

$ws.Range("H99").Value = 1001348.7
$ws.Range("I99").Value = 834689.5
$ws.Range("K99").Value = 834689.5
$ws.Range("M99").Value = -833191.5

$ws.Range("H122").Value = 3107.5715
$ws.Range("I122").Value = 1147.7273
$ws.Range("J122").Value = 10293.667
$ws.Range("K122").Value = 3443.1819
$ws.Range("L122").Value = 30881.001
$ws.Range("M122").Value = -993.1819
$ws.Range("N122").Value = -35781.001

$ws.Range("H126").Value = 1001348.7
$ws.Range("I126").Value = 834689.5
$ws.Range("K126").Value = 2504068.5
$ws.Range("M126").Value = -2501598.5

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 562380.3
$ws.Range("J141").Value = 1288887.4
$ws.Range("L141").Value = 1288887.4
$ws.Range("N141").Value = -1299247.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 349.0909
$ws.Range("I2").Value = 235.28572
$ws.Range("K2").Value = 1411.71432
$ws.Range("M2").Value = -1298.71432

$ws.Range("H4").Value = 2231.7073
$ws.Range("I4").Value = 1652.8572
$ws.Range("J4").Value = 2839.5
$ws.Range("K4").Value = 4958.571599999999
$ws.Range("L4").Value = 8518.5
$ws.Range("M4").Value = -4846.571599999999
$ws.Range("N4").Value = -8742.5

$ws.Range("H12").Value = 7.5
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 14
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 42
$ws.Range("M12").Value = 170
$ws.Range("N12").Value = -388

$ws.Range("H15").Value = 1020
$ws.Range("I15").Value = 40
$ws.Range("K15").Value = 120
$ws.Range("M15").Value = 20

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H39").Value = 3500
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 3500
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 10500
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -11088

$ws.Range("H55").Value = 4138.5483
$ws.Range("J55").Value = 4196.5576
$ws.Range("L55").Value = 12589.6728
$ws.Range("N55").Value = -12943.6728

$ws.Range("H68").Value = 1992
$ws.Range("I68").Value = 2242.25
$ws.Range("J68").Value = 1741.75
$ws.Range("K68").Value = 6726.75
$ws.Range("L68").Value = 5225.25
$ws.Range("M68").Value = -5915.75
$ws.Range("N68").Value = -6847.25

$ws.Range("H69").Value = 1467.6666
$ws.Range("J69").Value = 1499.4445
$ws.Range("L69").Value = 4498.333500000001
$ws.Range("N69").Value = -6120.333500000001

$ws.Range("H71").Value = 1992
$ws.Range("I71").Value = 2242.25
$ws.Range("J71").Value = 1741.75
$ws.Range("K71").Value = 20180.25
$ws.Range("L71").Value = 15675.75
$ws.Range("M71").Value = -16124.25
$ws.Range("N71").Value = -23787.75

$ws.Range("H72").Value = 1467.6666
$ws.Range("J72").Value = 1499.4445
$ws.Range("L72").Value = 13495.0005
$ws.Range("N72").Value = -21607.0005

$ws.Range("H74").Value = 7000
$ws.Range("J74").Value = 7000
$ws.Range("L74").Value = 21000
$ws.Range("N74").Value = -23122

$ws.Range("H77").Value = 7000
$ws.Range("J77").Value = 7000
$ws.Range("L77").Value = 63000
$ws.Range("N77").Value = -73608

$ws.Range("H94").Value = 449.75
$ws.Range("I94").Value = 449.75
$ws.Range("K94").Value = 1349.25
$ws.Range("M94").Value = -673.25

$ws.Range("H106").Value = 2091.6667
$ws.Range("J106").Value = 2091.6667
$ws.Range("L106").Value = 6275.000100000001
$ws.Range("N106").Value = -8167.000100000001

$ws.Range("H108").Value = 378.25
$ws.Range("I108").Value = 378.25
$ws.Range("K108").Value = 1134.75
$ws.Range("M108").Value = 1745.25

$ws.Range("H114").Value = 1410.6
$ws.Range("I114").Value = 1138.25
$ws.Range("J114").Value = 2500
$ws.Range("K114").Value = 3414.75
$ws.Range("L114").Value = 7500
$ws.Range("M114").Value = -160.75
$ws.Range("N114").Value = -14008

$ws.Range("H128").Value = 629995
$ws.Range("I128").Value = 629995
$ws.Range("K128").Value = 1889985
$ws.Range("M128").Value = -1885005

$ws.Range("H134").Value = 11666.667
$ws.Range("I134").Value = 5000
$ws.Range("K134").Value = 15000
$ws.Range("M134").Value = -9930

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 62999.668
$ws.Range("I94").Value = 61999
$ws.Range("K94").Value = 61999
$ws.Range("M94").Value = -61323

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H134").Value = 40326
$ws.Range("J134").Value = 40326
$ws.Range("L134").Value = 120978
$ws.Range("N134").Value = -126048

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8199.667
$ws.Range("I7").Value = 8199.667
$ws.Range("K7").Value = 8199.667
$ws.Range("M7").Value = -8087.666999999999

$ws.Range("H40").Value = 4000
$ws.Range("I40").Value = 4000
$ws.Range("K40").Value = 4000
$ws.Range("M40").Value = -3864

$ws.Range("H46").Value = 1462.5
$ws.Range("I46").Value = 1200
$ws.Range("J46").Value = 1725
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 1725
$ws.Range("M46").Value = -1012
$ws.Range("N46").Value = -2101

$ws.Range("H55").Value = 1045.75
$ws.Range("J55").Value = 1556.25
$ws.Range("L55").Value = 1556.25
$ws.Range("N55").Value = -1902.25

$ws.Range("H93").Value = 37037810
$ws.Range("I93").Value = 41667412
$ws.Range("K93").Value = 41667412
$ws.Range("M93").Value = -41666164

$ws.Range("H122").Value = 3586.8076
$ws.Range("I122").Value = 2910.8572
$ws.Range("K122").Value = 8732.5716
$ws.Range("M122").Value = -6282.571599999999

$ws.Range("H126").Value = 8199.667
$ws.Range("I126").Value = 8199.667
$ws.Range("K126").Value = 24599.001
$ws.Range("M126").Value = -22129.001

$ws.Range("H135").Value = 124994
$ws.Range("J135").Value = 124994
$ws.Range("L135").Value = 124994
$ws.Range("N135").Value = -135134

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1004
$ws.Range("I14").Value = 1004
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1004
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -836
$ws.Range("N14").ClearContents()

$ws.Range("H96").Value = 1850.75
$ws.Range("I96").Value = 1302.3334
$ws.Range("K96").Value = 1302.3334
$ws.Range("M96").Value = 70.66660000000002

$ws.Range("H107").Value = 285.14285
$ws.Range("I107").Value = 299.66666
$ws.Range("J107").Value = 198
$ws.Range("K107").Value = 898.9999799999999
$ws.Range("L107").Value = 594
$ws.Range("M107").Value = 1021.00002
$ws.Range("N107").Value = -4434

$ws.Range("H132").Value = 2716
$ws.Range("I132").Value = 2461.2
$ws.Range("J132").Value = 3990
$ws.Range("K132").Value = 7383.599999999999
$ws.Range("L132").Value = 11970
$ws.Range("M132").Value = -4853.599999999999
$ws.Range("N132").Value = -17030

$ws.Range("H140").Value = 35000
$ws.Range("J140").Value = 55000
$ws.Range("L140").Value = 55000
$ws.Range("N140").Value = -65360
